$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.579.35'
$ws.Range('E2').Value = '  +3.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.688.70'
$ws.Range('E3').Value = '  +8.60%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.88'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.680.92'
$ws.Range('E7').Value = '  +8.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('E8').Value = '  +4.84%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('E11').Value = '  +4.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.84'
$ws.Range('E12').Value = '  +3.25%  '
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.284.68'
$ws.Range('E14').Value = '  +8.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '684.02'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.697.92'
$ws.Range('E17').Value = '  +9.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '71.692.78'
$ws.Range('E18').Value = '  +3.34%  '
$ws.Range('E19').Value = '  +2.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.08'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.941'
$ws.Range('E22').Value = '  +3.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.27'
$ws.Range('E23').Value = '  +16.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.81'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '103.94'
$ws.Range('E25').Value = '  +2.68%  '
$ws.Range('E26').Value = '  +3.89%  '
$ws.Range('E27').Value = '  +5.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.18'
$ws.Range('E28').Value = '  +4.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.42'
$ws.Range('E29').Value = '  +5.89%  '
$ws.Range('E30').Value = '  +5.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.34'
$ws.Range('E31').Value = '  +6.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.21'
$ws.Range('E32').Value = '  +11.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '573.75'
$ws.Range('E33').Value = '  +3.33%  '
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  +3.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.41'
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.771.12'
$ws.Range('E37').Value = '  +4.61%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.146'
$ws.Range('E39').Value = '  +3.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0774'
$ws.Range('E40').Value = '  +3.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '35.36'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.78'
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0460'
$ws.Range('E44').Value = '  +8.07%  '
$ws.Range('E45').Value = '  +4.97%  '
$ws.Range('E46').Value = '  +7.88%  '
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('E48').Value = '  +4.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.43'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '134.46'
$ws.Range('E51').Value = '  +2.51%  '
